# TNRSequencer.xlsx edit:
#  - Log.addSTEPACTION / rename addSTEPSSGRP -> addSTEPBLOCK
#  - "LIST" sheet: A2 becomes the new test-case id "ZZ.010" (new shared
#    string), B2 gets the repetition count (5), A3 is cleared back to an
#    empty (but still styled) cell, and the old A4/A5 rows are removed.
#  - Update the active-cell selections on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LIST")
$ws2 = $wb.Worksheets.Item("Feuil1")

# A2: new shared string "ZZ.010" ; B2: repetition count 5
$ws1.Range("A2").Value = "ZZ.010"
$ws1.Range("B2").Value = 5

# A3 keeps its style but loses its text content
$ws1.Range("A3").ClearContents() | Out-Null

# A4 and A5 are removed outright (content + style), no row shifting
$ws1.Range("A4:A5").Clear() | Out-Null

# Update the selection on Feuil1 first ...
$ws2.Range("F5").Select() | Out-Null

# ... then select on LIST last so it remains the active/selected tab
$ws1.Range("B2").Select() | Out-Null
